$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank rows at 18:19 -- shifts the SUM/Total rows (and the
#    G:H merge) down from 21/22 to 23/24, same as the authored diff.
$ws.Range("A18:A19").EntireRow.Insert()

# 2. Add the purchase-date entry for the existing "solder paste" row (16).
#    Register the "mm-dd-yy" (builtin numFmtId 14) date format first so the
#    the border style ends up before it in cellXfs, matching target order.
$ws.Rows.Item(22).Borders.Item(9).LineStyle = 1

$ws.Range("C16").NumberFormat = "mm-dd-yy"
$ws.Range("C16").Value = "5/12/2023"

# 3. New rows 18 & 19: contour tool / electronic silicone purchases.
$ws.Range("A18").Value = "contour tool"
$ws.Range("C18").NumberFormat = "mm-dd-yy"
$ws.Range("C18").Value = "4/19/2023"
$ws.Range("H18").Value = 11.79

$ws.Range("A19").Value = "electronic silicone"
$ws.Range("C19").NumberFormat = "mm-dd-yy"
$ws.Range("C19").Value = "8/7/2023"
$ws.Range("H19").Value = 10.67

# 4. Pump/tray tray-dry totals block (rows 39-43).
$ws.Range("H39").Value = 3.98
$ws.Range("H40").Value = 6.98
$ws.Range("H41").Value = 10.48
$ws.Range("H42").Value = 8.98
$ws.Range("H43").Formula = "=H23+SUM(H39:H42)"
